# Update extrapolation calibration values to remove noisy sub-$5 price rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 113573.8422524486
$ws.Range("E4").Value = -0.02389311441453535
$ws.Range("F4").Value = 0.1863027698936319
$ws.Range("G4").Value = -1.420897712649651
$ws.Range("H4").Value = 12.30725302030576

$ws.Range("D6").Value = 114907.6258735845
$ws.Range("E6").Value = -0.030215779576078
$ws.Range("F6").Value = 0.2297571791027435
$ws.Range("G6").Value = -1.475677235959236
$ws.Range("H6").Value = 11.46325218886858

$ws.Range("D8").Value = 116648.0593894568
$ws.Range("E8").Value = -0.044540681786888
$ws.Range("F8").Value = 0.2047575199294869
$ws.Range("G8").Value = -1.069679848952052
$ws.Range("H8").Value = 7.442129706015947

$ws.Range("D9").Value = 118197.9349878586
$ws.Range("E9").Value = -0.07141440347529397
$ws.Range("F9").Value = 0.3179644234804466
$ws.Range("G9").Value = -1.688448949661576
$ws.Range("H9").Value = 10.77396189580011

$ws.Range("D10").Value = 119556.9225141097
$ws.Range("E10").Value = -0.1174893032825472
$ws.Range("F10").Value = 0.4412331378602905
$ws.Range("G10").Value = -1.88552703207416
$ws.Range("H10").Value = 9.522205608575337

$ws.Range("D15").Value = 112756.5884329135
$ws.Range("E15").Value = -0.02350322210009533
$ws.Range("F15").Value = 0.1327488096091299
$ws.Range("G15").Value = -0.6607562226556706
$ws.Range("H15").Value = 7.546967890305524

$ws.Range("D18").Value = 112807.5780221454
$ws.Range("E18").Value = -0.02665615810723996
$ws.Range("F18").Value = 0.1628785304374701
$ws.Range("G18").Value = -0.6252628023087783
$ws.Range("H18").Value = 6.294792007321606
